$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new rows 106-129 (pendientes) ---
$ws.Range("A106").Value = "AL RECREAR INGRESAR FECHA VALIDA ACTUALMENTE LA DEJA CON FORMATO DE INICIO, VALIDAR COLONIA, NOeXSTERIOR MENOR A 3"
$ws.Range("B106").Value = "Pendiente"
$ws.Rows.Item(106).RowHeight = 45

$ws.Range("A108").Value = "Deshabilitar el agregar a las imagens de extras"
$ws.Range("B108").Value = "Pendiente"

$ws.Range("A109").Value = "REVISAR TOKEN LA EXPIRACIÓN"
$ws.Range("B109").Value = "Pendiente"

$ws.Range("A110").Value = "Ocr desahabilitar visivilidad"
$ws.Range("B110").Value = "Pendiente"

$ws.Range("A111").Value = "revisar los success de las pantallas de nueva solicitud en los inputs llenados por medio del OCR"
$ws.Range("B111").Value = "Pendiente"
$ws.Rows.Item(111).RowHeight = 30

$ws.Range("A113").Value = "validar el contenido de los registros regresador por OCR y mostrar o no dependiendo del mismo, quitar las validaciones para eventos de inicio"
$ws.Range("B113").Value = "Pendiente"
$ws.Range("C113").Value = "las validaciones d elos campos se dispararán al dar clic sobre los mismos, cambios, blur y al dar clic en siguiente"
$ws.Rows.Item(113).RowHeight = 45

$ws.Range("A114").Value = "Resetear los controles y validaciones al dar click hacia atrás en header futeer"
$ws.Range("B114").Value = "Pendiente"
$ws.Rows.Item(114).RowHeight = 30

$ws.Range("A115").Value = "quitar colores a las etiquetas de los validator"
$ws.Range("B115").Value = "Pendiente"

$ws.Range("A116").Value = "habilitar spiner con un max de 2"
$ws.Range("B116").Value = "Pendiente"

$ws.Range("A117").Value = "revision de correo,  pendejo esto ya estaba "
$ws.Range("B117").Value = "Pendiente"

$ws.Range("A118").Value = "validar los telefonos colocar paloma hasta 10 digitos menos de estos colocar amarillo"
$ws.Range("B118").Value = "Pendiente"
$ws.Rows.Item(118).RowHeight = 30

$ws.Range("A119").Value = "que no se borre el numero de telefono si es menor que 10 en el blur, que no mame"
$ws.Range("B119").Value = "Pendiente"
$ws.Rows.Item(119).RowHeight = 30

$ws.Range("A120").Value = "hacer Uppercase en todas las cajas de texto "
$ws.Range("B120").Value = "Pendiente"

$ws.Range("A121").Value = "COLOCAR VERDE EN EL HEADER  AL VALIDAR TODO SUCCESS EN EL CLICK ABAJO"
$ws.Range("B121").Value = "Pendiente"
$ws.Rows.Item(121).RowHeight = 30

$ws.Range("A122").Value = "EN LA PERSONA POLITICA QUITAR CONTENIDO DE inputs si se eleje si y no despues"
$ws.Range("B122").Value = "Pendiente"
$ws.Rows.Item(122).RowHeight = 30

$ws.Range("A123").Value = "para las referencias familiares validar el llenado de los mismo al meter cualquier campo"
$ws.Range("B123").Value = "Pendiente"
$ws.Rows.Item(123).RowHeight = 30

$ws.Range("A124").Value = "acomodar el pad de firma"
$ws.Range("B124").Value = "Pendiente"

$ws.Range("A125").Value = "validar las imágenes de INE frente atrás y FIRMA como minimo"
$ws.Range("B125").Value = "Pendiente"
$ws.Rows.Item(125).RowHeight = 30

$ws.Range("A126").Value = "revisar el envio de la solicitud"
$ws.Range("B126").Value = "Pendiente"

$ws.Range("A127").Value = "GenerarBitacora de Operaciones"
$ws.Range("B127").Value = "Pendiente"
$ws.Range("A127:C127").Interior.Color = 65535

$ws.Range("A128").Value = "modificar la imagen en la presentacion de documentos para que se vea completa"
$ws.Range("B128").Value = "Pendiente"
$ws.Rows.Item(128).RowHeight = 30

$ws.Range("A129").Value = "Validar si el usuario no ingresa la informacion, click por el header"
$ws.Rows.Item(129).RowHeight = 30

